$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# --- Fix row 19: "LIBRARY SOURCE" display name was mapped to the wrong
# node/property (it was a duplicate of row 18 "genomic_info.library_layout").
# Point it at the correct genomic_info.library_source field; the TEXTSPLIT
# array formula in E19 (spilling into F19) recalculates automatically.
$ws.Range("D19").Value = "genomic_info.library_source"

# --- Append a new mapping row for "Study Access" (row 84), following the
# same Page/Area/Display Name/Full Name/Node Name/Property Name pattern as
# the existing rows.
$ws.Range("A84").Value = "Study Access"
$ws.Range("B84").Value = "Study Access"
$ws.Range("C84").Value = "Study"
$ws.Range("D84").Value = "study.study_access"
$ws.Range("E84").Value = "study"
$ws.Range("F84").Value = "study_access"

# Move the active selection to reflect where the editor ended up working.
[void]$ws.Range("E9").Select()
